$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 38, shifting existing rows 38-176 down to 39-177.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly data entry.
$ws.Cells.Item(38, 1).Value = 3
$ws.Cells.Item(38, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44676
$ws.Cells.Item(38, 5).Value = 5
$ws.Cells.Item(38, 6).Value = 100112052
$ws.Cells.Item(38, 7).Value = "Albahaca"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 60
$ws.Cells.Item(38, 11).Value = 4500
$ws.Cells.Item(38, 12).Value = 4500
$ws.Cells.Item(38, 13).Value = 4500
$ws.Cells.Item(38, 14).Value = "$/docena de matas"
$ws.Cells.Item(38, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(38, 16).Value = 750
$ws.Cells.Item(38, 17).Value = 6
$ws.Cells.Item(38, 18).Value = "Hortaliza"
